$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("O_TransactionActivity")

for ($r = 3; $r -le 24; $r++) {
    $cell = $ws.Cells.Item($r, 6)  # Column F
    if ($cell.Value() -eq "LOAN1") {
        $cell.Value = "IDHJ-EGNY"
    }
}
